$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2) -----------------------------------------------
# Populated in this exact order so the shared-string table comes out in
# the same order the source workbook has it in.
$ws.Range("C2").Value = "Books"
$ws.Range("D2").Value = "Page"
$ws.Range("E2").Value = "Reference content"
$ws.Range("F2").Value = "Link"
$ws.Range("B2").Value = "Index"

# --- Row 3 --------------------------------------------------------------
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Applied Cryptography : Protocols, Algorithms and Source Code in C"
$ws.Range("D3").Value = 508
$ws.Range("E3").Value = "Foundations and applications of LFSR system"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://drm.phenikaa-uni.edu.vn/vi/handle/DRM/6320?status=activate", "", "", "https://drm.phenikaa-uni.edu.vn/vi/handle/DRM/6320?status=activate")
$ws.Range("F3").Value = "Digital Right Management"

# --- Row 4 --------------------------------------------------------------
$ws.Range("B4").Value = 2
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.01signal.com/other/lfsr-galois-fibonacci/", "", "", "https://www.01signal.com/other/lfsr-galois-fibonacci/")
$ws.Range("F4").Value = "01signal: Conversion between Galois and Fibonacci polynomials of Linear-Feedback Shift Register"
$ws.Range("C4").Value = "website: 01 signal"
$ws.Range("D4").Value = ".. "
$ws.Range("E4").Value = "Expression of LFSR in fibonacci and galois field"

# --- Index continuation, rows 5-10 --------------------------------------
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 6
$ws.Range("B9").Value = 7
$ws.Range("B10").Value = 8

# --- Column widths ---------------------------------------------------------
# Target widths (24.21875 / 8.88671875 / 25.77734375 chars) are finer-grained
# than this host's column-width quantization (1/6-character steps), so these
# inputs are chosen to land on the closest achievable stored width.
$ws.Range("C1").ColumnWidth = 23.3337
$ws.Range("D1").ColumnWidth = 8.0017
$ws.Range("E1").ColumnWidth = 25.0023

# --- View: selection + zoom ----------------------------------------------
$ws.Range("C5").Select() | Out-Null
$excel.ActiveWindow.Zoom = 205 | Out-Null
